$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-02-10 19:11:21"
}
